$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (copy formatting from H1, then set text)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-60
$data = @(
    @(2, 7, 8),
    @(3, 6, 6),
    @(4, 6, 6),
    @(5, 6, 6),
    @(6, 6, 6),
    @(7, 7, 7),
    @(8, 6, 6),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 9, 9),
    @(12, 6, 6),
    @(13, 10, 10),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 6, 6),
    @(17, 10, 10),
    @(18, 8, 8),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 7, 7),
    @(33, 8, 8),
    @(34, 8, 8),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 7, 7),
    @(38, 7, 7),
    @(39, 8, 8),
    @(40, 7, 8),
    @(41, 8, 9),
    @(42, 7, 7),
    @(43, 8, 8),
    @(44, 8, 8),
    @(45, 5, 6),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 7, 7),
    @(50, 7, 7),
    @(51, 8, 8),
    @(52, 7, 7),
    @(53, 9, 9),
    @(54, 7, 7),
    @(55, 7, 7),
    @(56, 7, 7),
    @(57, 4, 4),
    @(58, 4, 4),
    @(59, 4, 4),
    @(60, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
